# Correção de erro excel
# - D3: "l" -> "K"
# - I3: "123" -> "987654321"
# - T3: "Funciona ;)" -> "Avariado"

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Folha1")

$ws.Range("D3").Value = "K"

$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "987654321"

$ws.Range("T3").Value = "Avariado"
